# "COal 2030 cap, fix UCs and WH Elc growth"
#
# On the "RSD" worksheet the standalone FT-RSDKER (kerosene) UC_ACT cap
# table (rows 7-9) is removed; the FT-RSDCOA (coal) cap table that used
# to sit right underneath it (rows 10-12, with the 2030/2050 caps now
# tightened to 4/2) takes its place at the top, directly under a new
# "~TFM_INS" marker placed at B5. The remaining FT-RSDGAS table (with its
# own "~TFM_INS" header, previously starting at row 17) simply shifts up
# to start at row 14 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSD")

# Remove the old FT-RSDKER table (rows 7-9). This shifts all the
# following rows up by three: the coal (FT-RSDCOA) table that used to be
# rows 10-12 becomes the new rows 7-9, and the second "~TFM_INS" /
# FT-RSDGAS block (old rows 17-21) becomes rows 14-18 - exactly matching
# the target layout.
$ws.Rows("7:9").Delete()

# Add the new "~TFM_INS" marker above the (now relocated) coal table.
$ws.Range("B5").Value = "~TFM_INS"

# Update the saved selection to match the edited workbook.
$ws.Range("D12").Select()
